$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 319 (pushes the existing row 319.. down by one) ---
$ws.Rows.Item(319).EntireRow.Insert()
$ws.Cells.Item(319,1).Value = 4
$ws.Cells.Item(319,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(319,3).Value = "Los Lagos"
$ws.Cells.Item(319,4).Value = 45006
$ws.Cells.Item(319,5).Value = 10
$ws.Cells.Item(319,6).Value = 100114014
$ws.Cells.Item(319,7).Value = "Betarraga"
$ws.Cells.Item(319,8).Value = "Sin especificar"
$ws.Cells.Item(319,9).Value = "Primera"
$ws.Cells.Item(319,10).Value = 1000
$ws.Cells.Item(319,11).Value = 1100
$ws.Cells.Item(319,12).Value = 1200
$ws.Cells.Item(319,13).Value = 1150
$ws.Cells.Item(319,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(319,15).Value = "Provincia de Cautín"
$ws.Cells.Item(319,16).Value = 230
$ws.Cells.Item(319,17).Value = 5
$ws.Cells.Item(319,18).Value = "Hortaliza"

# --- Insert a second new record at row 371 (after the first insert has already
#     shifted everything down by one) ---
$ws.Rows.Item(371).EntireRow.Insert()
$ws.Cells.Item(371,1).Value = 4
$ws.Cells.Item(371,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(371,3).Value = "Los Lagos"
$ws.Cells.Item(371,4).Value = 45005
$ws.Cells.Item(371,5).Value = 10
$ws.Cells.Item(371,6).Value = 100114014
$ws.Cells.Item(371,7).Value = "Betarraga"
$ws.Cells.Item(371,8).Value = "Sin especificar"
$ws.Cells.Item(371,9).Value = "Primera"
$ws.Cells.Item(371,10).Value = 500
$ws.Cells.Item(371,11).Value = 1100
$ws.Cells.Item(371,12).Value = 1200
$ws.Cells.Item(371,13).Value = 1150
$ws.Cells.Item(371,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(371,15).Value = "Provincia de Cautín"
$ws.Cells.Item(371,16).Value = 230
$ws.Cells.Item(371,17).Value = 5
$ws.Cells.Item(371,18).Value = "Hortaliza"
